$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("curso"), shifting curso and the
# question columns one to the right (D:J -> E:K).
$ws.Columns.Item(4).Insert()

# New column D: "marca" (brand) header + values
$ws.Range("D1").Value = "marca"
$ws.Range("D3").Value = "maestro"
$ws.Range("D2").Value = "sodimac"
$ws.Range("D4").Value = "sodimac"

# Re-affirm curso column (shifted from D to E) stays intact
$ws.Range("E1").Value = "curso"
$ws.Range("E2").Value = "ferreteria"
$ws.Range("E3").Value = "pintura"
$ws.Range("E4").Value = "electricista"

# Re-affirm the quiz answer columns (shifted from E:J to F:K)
$ws.Range("F1").Value = "¿como a?"
$ws.Range("G1").Value = "¿que b?"
$ws.Range("H1").Value = "¿como c?"
$ws.Range("I1").Value = "¿cual d?"
$ws.Range("J1").Value = "¿como e?"
$ws.Range("K1").Value = "¿que f?"

$ws.Range("F2").Value = "a"
$ws.Range("G2").Value = "a"
$ws.Range("H2").Value = "c"
$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "a"
$ws.Range("K2").Value = "c"

$ws.Range("F3").Value = "a"
$ws.Range("G3").Value = "c"
$ws.Range("H3").Value = "b"
$ws.Range("I3").Value = "c"
$ws.Range("J3").Value = "a"
$ws.Range("K3").Value = "a"

$ws.Range("F4").Value = "a"
$ws.Range("G4").Value = "c"
$ws.Range("H4").Value = "b"
$ws.Range("I4").Value = "c"
$ws.Range("J4").Value = "a"
$ws.Range("K4").Value = "a"

# Restore the wrap-text style on F1 (shifted from E1, which had it before)
$ws.Range("F1").WrapText = $true

# Update the active selection to match the target view state
$ws.Range("H7").Select()
